$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and restore the Decentraland/EnergySwap row order)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.485.83"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.082.32"
$ws.Range("E3").Value = "  +4.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.03"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5216"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4332"
$ws.Range("E8").Value = "  +4.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08772"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.00"
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.156"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.27"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.087.64"
$ws.Range("E13").Value = "  +4.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.654"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.744"
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.84"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001119"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06635"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.79"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.251"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.580.41"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.32"
$ws.Range("E24").Value = "  +4.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.333.90"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.23"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.530"
$ws.Range("E28").Value = "  +5.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.37"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.37"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.187"
$ws.Range("E31").Value = "  +4.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.069"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.533"
$ws.Range("E34").Value = "  +14.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.832"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02567"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.675"
$ws.Range("E37").Value = "  +7.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.469"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06657"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2247"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.45"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6729"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6331"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.79"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.198"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.618"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.232"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.14"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.182"
$ws.Range("E51").Value = "  +6.86%  "

Write-Output "Updated crypto list values"
